# LIVEHTA-723 / LIVEHTA-931: bump LiveSLR build number shown in the
# "Application version data" sheet (row 2, column B) to Build #52677.
#
# (The workbook's fileVersion/rupBuild, xr:revisionPtr documentId, and the
# theme's thm15:themeFamily name are Office/authoring-tool bookkeeping
# metadata that Excel regenerates itself on save and are not part of the
# visible workbook content reachable through the object model.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = "Copyright @ 2023 Cytel Inc. LiveSLR 4.0.0.0 - Build #52677"
